# Update the "取得日時" (acquisition datetime) column (A) for rows 2-10
# on the first worksheet from "2025-12-07 01:29:42" to "2025-12-07 02:02:01",
# as described in the commit: "Append: 2025-12-07 02:02 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-07 02:02:01"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
